# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Uva" (Terminal La Palmera de La Serena)
# as row 81, pushing the existing rows 81-98 down to 82-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81 (existing rows 81..98 shift to 82..99,
# and Excel's "insert" copies the formatting from the row above, which already
# gives the date column D the expected numeric/date style).
$ws.Rows.Item(81).Insert()

$ws.Range("A81").Value2 = 8
$ws.Range("B81").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C81").Value2 = "Coquimbo"
$ws.Range("D81").Value2 = 44637
$ws.Range("E81").Value2 = 4
$ws.Range("F81").Value2 = "Fruta"
$ws.Range("G81").Value2 = 100109
$ws.Range("H81").Value2 = "Uva"
$ws.Range("I81").Value2 = 100109001
$ws.Range("J81").Value2 = "Uva"
$ws.Range("K81").Value2 = "Red Globe"
$ws.Range("L81").Value2 = "Primera"
$ws.Range("M81").Value2 = 400
$ws.Range("N81").Value2 = 9500
$ws.Range("O81").Value2 = 10000
$ws.Range("P81").Value2 = 9750
$ws.Range("Q81").Value2 = "$/bandeja 18 kilos"
$ws.Range("R81").Value2 = "Provincia del Elquí"
$ws.Range("S81").Value2 = 542
$ws.Range("T81").Value2 = 18
